# Fruta / hortaliza, semanal
# Insert a new daily price record as the new row 3, pushing all existing
# records (old rows 3..122) down by one row (to rows 4..123).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3 (shifts rows 3:122 -> 4:123)
$ws.Rows(3).Insert()

# Populate the newly inserted row 3 with the new record's data
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 45160
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100108
$ws.Range("H3").Value = "Tropicales y subtropicales"
$ws.Range("I3").Value = 100108007
$ws.Range("J3").Value = "Coco"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 36000
$ws.Range("O3").Value = 36000
$ws.Range("P3").Value = 36000
$ws.Range("Q3").Value = "`$/malla 20 unidades"
$ws.Range("R3").Value = "Perú"
$ws.Range("S3").Value = 1800
$ws.Range("T3").Value = 20
